$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "severity" column (C), styled like the existing
# header cells (A1/B1 use style index 1 - bold, bordered, centered).
# Copy A1's formatting onto C1 first (Copy carries the cell style along),
# then overwrite the copied text with the real header label.
$ws.Range("A1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "severity"

# Row 2 ("Ahmoq") gets severity 55; every other data row (3-157) gets
# the default severity of 1.
$ws.Range("C2").Value = 55

for ($r = 3; $r -le 157; $r++) {
    $ws.Cells.Item($r, 3).Value = 1
}
